$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldCommitText = "IndicatorQuantiles.R, Git Commit ID: 54e4488a188edf59eafc4b9cfe53dc7125db7b32"
$newCommitText = "IndicatorQuantiles.R, Git Commit ID: 0e4152332be22faf035a2e2fc83ad2cca4c8a7fc"

$oldPid = 30656
$newPid = 25596

for ($r = 2; $r -le 80; $r++) {
    # Column AJ (36) holds the "IndicatorQuantiles.R, Git Commit ID: ..." text
    $ajCell = $ws.Cells.Item($r, 36)
    $ajValue = $ajCell.Value()
    if ($ajValue -eq $oldCommitText) {
        $ajCell.Value = $newCommitText
    }

    # Column AH (34) holds the "pid" numeric value
    $ahCell = $ws.Cells.Item($r, 34)
    $ahValue = $ahCell.Value()
    if ($ahValue -eq $oldPid) {
        $ahCell.Value = $newPid
    }
}
